$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The NATMI TPM post-processing script was re-run with updated normalization.
# The recomputed output table no longer reports a "Target cluster = MuSCs"
# breakdown, so those three rows (old rows 8-10) are removed, and every
# remaining data row (old rows 2-7) is overwritten with the freshly computed
# expression / specificity statistics.
$ws.Rows("8:10").Delete()

$data = New-Object 'object[,]' 6,20
$data[0,0] = "ECs"
$data[0,1] = "Ccl2"
$data[0,2] = "Ccr5"
$data[0,3] = "ECs"
$data[0,4] = 2
$data[0,5] = 0.6666666666666666
$data[0,6] = 0.07794266666666667
$data[0,7] = 0.233828
$data[0,8] = 0.002827880818927331
$data[0,9] = 0.00282788081892733
$data[0,10] = 1
$data[0,11] = 0.3333333333333333
$data[0,12] = 0.001937666666666667
$data[0,13] = 0.005813
$data[0,14] = 0.0230007399171451
$data[0,15] = 0.02300073991714511
$data[0,16] = 0.0001510269071111111
$data[0,17] = 0.001359242164
$data[0,18] = [double]"6.504335123283085E-05"
$data[0,19] = [double]"6.504335123283083E-05"
$data[1,0] = "ECs"
$data[1,1] = "Ccl2"
$data[1,2] = "Ccr5"
$data[1,3] = "FAPs"
$data[1,4] = 2
$data[1,5] = 0.6666666666666666
$data[1,6] = 0.07794266666666667
$data[1,7] = 0.233828
$data[1,8] = 0.002827880818927331
$data[1,9] = 0.00282788081892733
$data[1,10] = 1
$data[1,11] = 0.3333333333333333
$data[1,12] = 0.082306
$data[1,13] = 0.246918
$data[1,14] = 0.9769992600828549
$data[1,15] = 0.976999260082855
$data[1,16] = 0.006415149122666668
$data[1,17] = 0.057736342104
$data[1,18] = 0.0027628374676945
$data[1,19] = 0.0027628374676945
$data[2,0] = "FAPs"
$data[2,1] = "Ccl2"
$data[2,2] = "Ccr5"
$data[2,3] = "ECs"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 25.22419466666667
$data[2,7] = 75.672584
$data[2,8] = 0.9151728997907317
$data[2,9] = 0.9151728997907316
$data[2,10] = 1
$data[2,11] = 0.3333333333333333
$data[2,12] = 0.001937666666666667
$data[2,13] = 0.005813
$data[2,14] = 0.0230007399171451
$data[2,15] = 0.02300073991714511
$data[2,16] = 0.0488760811991111
$data[2,17] = 0.439884730792
$data[2,18] = 0.02104965384730612
$data[2,19] = 0.02104965384730612
$data[3,0] = "FAPs"
$data[3,1] = "Ccl2"
$data[3,2] = "Ccr5"
$data[3,3] = "FAPs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 25.22419466666667
$data[3,7] = 75.672584
$data[3,8] = 0.9151728997907317
$data[3,9] = 0.9151728997907316
$data[3,10] = 1
$data[3,11] = 0.3333333333333333
$data[3,12] = 0.082306
$data[3,13] = 0.246918
$data[3,14] = 0.9769992600828549
$data[3,15] = 0.976999260082855
$data[3,16] = 2.076102566234667
$data[3,17] = 18.684923096112
$data[3,18] = 0.8941232459434255
$data[3,19] = 0.8941232459434255
$data[4,0] = "MuSCs"
$data[4,1] = "Ccl2"
$data[4,2] = "Ccr5"
$data[4,3] = "ECs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 2.260080333333333
$data[4,7] = 6.780241
$data[4,8] = 0.08199921939034102
$data[4,9] = 0.08199921939034102
$data[4,10] = 1
$data[4,11] = 0.3333333333333333
$data[4,12] = 0.001937666666666667
$data[4,13] = 0.005813
$data[4,14] = 0.0230007399171451
$data[4,15] = 0.02300073991714511
$data[4,16] = 0.004379282325888889
$data[4,17] = 0.039413540933
$data[4,18] = 0.001886042718606155
$data[4,19] = 0.001886042718606156
$data[5,0] = "MuSCs"
$data[5,1] = "Ccl2"
$data[5,2] = "Ccr5"
$data[5,3] = "FAPs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 2.260080333333333
$data[5,7] = 6.780241
$data[5,8] = 0.08199921939034102
$data[5,9] = 0.08199921939034102
$data[5,10] = 1
$data[5,11] = 0.3333333333333333
$data[5,12] = 0.082306
$data[5,13] = 0.246918
$data[5,14] = 0.9769992600828549
$data[5,15] = 0.976999260082855
$data[5,16] = 0.1860181719153333
$data[5,17] = 1.674163547238
$data[5,18] = 0.08011317667173487
$data[5,19] = 0.08011317667173487

$ws.Range("A2:T7").Value = $data
